$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title / byline / email ---
Replace-Text "Tranquility Found: Musings on Inner Peace" "Exploring the Heart of Democracy: Government and Its Role in Shaping Our World"
Replace-Text "Clarissa Shaw" "Alex Roberts"
Replace-Text "ClarissaShaw@rhymemail" "alex"
Replace-Text "com" "roberts@validedu.org"

# --- Body paragraph 1 ---
Replace-Text "We are all chasing after something, some yearning that we call happiness, fulfillment, contentment, or peace" "In the tapestry of human existence, nestled amidst the myriad threads of life's intricacies, lies the concept of government: a construct that weaves together the fabric of societies and shapes our collective destinies"
Replace-Text " While these concepts often seem elusive and unattainable, they lie within us, waiting to be discovered" " Government, a fascinating entity that has intrigued scholars, philosophers, and citizens alike throughout history, serves as the cornerstone of organized societies"
Replace-Text " In a world of ceaseless noise, it is essential to seek tranquility and find solace in the sanctuary of our own minds" " It is the engine that drives progress, the guardian of rights, and the arbiter of justice"

Replace-Text " Like the gentle lapping of waves on a tranquil shore, inner peace is a state of calmness, free from agitation and disturbance. It is a journey that begins with self-reflection, where we unravel the tangled threads of our emotions, thoughts, and actions to understand ourselves better. From this place of self-discovery, we can cultivate resilience, gratitude, compassion, and purpose, the pillars upon which inner peace is built" " Embarking on a journey to unravel the complexities of government and its profound influence on our world, we will explore its fundamental principles, its diverse forms, and the intricate interplay between citizens and their governing bodies"

Replace-Text "Through the exploration of our inner landscapes, we can uncover the source of our disquiet and work towards purging them" "In the realm of politics, we delve into the dynamic interplay of power, ideology, and public opinion"
Replace-Text " With every acknowledgment of our anxieties, fears, and desires, we slowly chip away at the barriers we have erected around our hearts" " We uncover the mechanisms through which decisions are made, policies are formulated, and leaders are chosen"
Replace-Text " Forgiveness, both of ourselves and others, becomes a fundamental step in this process" " Analyzing historical and contemporary political landscapes, we grapple with questions of representation, accountability, and the delicate balance between individual liberties and collective well-being"

Replace-Text " We learn to let go of grudges, resentments, and wounds, choosing instead to embrace acceptance and understanding. Like a gentle breeze sweeping away dust and cobwebs, forgiveness clears the path to tranquility, allowing us to move forward with lightness and grace" " We trace the evolution of political thought, examining the ideas of influential thinkers and movements that have shaped our understanding of governance"

Replace-Text "To find inner peace, we must cultivate gratitude and appreciation for the myriad blessings in our lives" "Delving into the annals of history, we embark on a voyage to explore the diverse forms that governments have taken across time and space"
Replace-Text " Like a child marveling at the beauty of a blooming flower, we must train our eyes to see the wonders that surround us" " From ancient civilizations to modern nation-states, we uncover the unique characteristics, challenges, and achievements of different political systems"
Replace-Text " Each breath we take, each sunrise we witness, each act of kindness we receive is an opportunity to find gratitude" " We examine the rise and fall of empires, the birth of democracies, and the struggles for independence that have shaped the world's political map"

Replace-Text " By acknowledging the abundance in our lives, we open ourselves to a profound sense of contentment and inner peace. Our hearts soften, like clay in the hands of a potter, as we learn to appreciate the simple joys of existence" " Through this journey, we gain a deeper appreciation for the complexity of governance and the diverse approaches to organizing and managing human societies"

# --- Summary paragraph ---
Replace-Text "Inner peace is a sanctuary we can access within ourselves, a place of tranquility where the storms of the world cannot reach" "In exploring the intricate world of government, we gain insights into the mechanisms that shape our collective lives"
Replace-Text " Through self-reflection, we can identify and release the sources of our discontent" " We delve into the realm of politics, examining the dynamics of power, "
Replace-Text " Forgiveness and gratitude are the keys that unlock the doors to inner " ""
Replace-Text "peace, allowing us to move forward with lightness and appreciation" "ideology, and public opinion that drive decision-making and policy formulation. We trace the evolution of political thought and its impact on governance"
Replace-Text " By cultivating these virtues and connecting with the abundance of life, we can discover the tranquility that lies at the heart of our being" " Additionally, we delve into the diverse forms governments take across time and space, appreciating the unique challenges and achievements of different political systems. Ultimately, this exploration serves to deepen our understanding and appreciation for the vital role government plays in shaping our world"

# --- New trailing empty paragraph ---
$d.Content.InsertParagraphAfter()
